$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed values.
# D-column values are plain text (not locale-safe numbers, e.g. "63.479.40"),
# so force text formatting before assigning to avoid Excel auto-converting them
# to numbers (which would also strip meaningful trailing zeros, e.g. "4.30" -> 4.3).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.479.40'
$ws.Range("E2").Value = '  +4.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.469.32'
$ws.Range("E3").Value = '  +5.88%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.86'
$ws.Range("E5").Value = '  +3.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.44'
$ws.Range("E6").Value = '  +9.55%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  +2.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.468.09'
$ws.Range("E9").Value = '  +5.82%  '
$ws.Range("E10").Value = '  +4.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.73'
$ws.Range("E11").Value = '  +3.73%  '
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +4.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.37'
$ws.Range("E14").Value = '  +11.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.911.96'
$ws.Range("E15").Value = '  +6.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.349.39'
$ws.Range("E16").Value = '  +4.57%  '
$ws.Range("E17").Value = '  +6.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.471.28'
$ws.Range("E18").Value = '  +6.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.27'
$ws.Range("E19").Value = '  +5.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.41'
$ws.Range("E20").Value = '  +8.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("E21").Value = '  +5.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.82'
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.74'
$ws.Range("E24").Value = '  +2.58%  '
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  +8.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.17'
$ws.Range("E28").Value = '  +3.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.34'
$ws.Range("E29").Value = '  +7.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0819'
$ws.Range("E30").Value = '  +12.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.83'
$ws.Range("E31").Value = '  +14.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.85'
$ws.Range("E32").Value = '  +7.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '176.57'
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("E34").Value = '  +10.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.400'
$ws.Range("E35").Value = '  +4.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.94'
$ws.Range("E36").Value = '  +5.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '373.29'
$ws.Range("E37").Value = '  +15.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.46'
$ws.Range("E38").Value = '  +8.08%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  +12.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.48'
$ws.Range("E42").Value = '  +6.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '150.93'
$ws.Range("E43").Value = '  +9.88%  '
$ws.Range("E44").Value = '  +5.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.77'
$ws.Range("E45").Value = '  +8.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.600'
$ws.Range("E46").Value = '  +5.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0965'
$ws.Range("E47").Value = '  +2.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0523'
$ws.Range("E48").Value = '  +4.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0239'
$ws.Range("E49").Value = '  +9.17%  '
$ws.Range("E50").Value = '  +4.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.15'
$ws.Range("E51").Value = '  +6.97%  '
